$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J, matching the style of the existing
# header cell H1 (bold / centered / bordered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-21 for columns I (I0) and J (IF)
$data = @(
    @(2, 2, 2),
    @(3, 8, 8),
    @(4, 8, 8),
    @(5, 5, 5),
    @(6, 9, 9),
    @(7, 7, 7),
    @(8, 8, 8),
    @(9, 7, 8),
    @(10, 9, 9),
    @(11, 8, 8),
    @(12, 8, 8),
    @(13, 6, 6),
    @(14, 9, 9),
    @(15, 9, 9),
    @(16, 7, 8),
    @(17, 7, 8),
    @(18, 6, 8),
    @(19, 5, 6),
    @(20, 1, 1),
    @(21, 9, 9)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
